$wb = $excel.ActiveWorkbook

# Rename the "Requested quantity" header on the existing sheets.
$wsWeekly = $wb.Worksheets.Item(1)
$wsWeekly.Range("B1").Value = "Weekly_PO_Qty"

$wsMonthly = $wb.Worksheets.Item(2)
$wsMonthly.Range("B1").Value = "Monthly_PO_Qty"

# Add a new "PO Forecast" sheet at the end of the workbook.
$sheetCount = $wb.Worksheets.Count
$wsForecast = $wb.Worksheets.Add($null, $wb.Worksheets.Item($sheetCount))
$wsForecast.Name = "PO Forecast"

# Header row.
$wsForecast.Range("A1").Value = "ds"
$wsForecast.Range("B1").Value = "PO_Forecast"
$wsForecast.Range("C1").Value = "yhat_lower"
$wsForecast.Range("D1").Value = "yhat_upper"

# Match the header formatting (bold, centered, bordered) used on the other sheets.
$wsWeekly.Range("A1:B1").Copy()
$wsForecast.Range("A1:D1").PasteSpecial(-4122)

# Forecast data rows (ds, PO_Forecast, yhat_lower, yhat_upper).
$data = New-Object 'object[,]' 60,4
$data[0,0] = 45081.99999999999
$data[0,1] = 152
$data[0,2] = -359.810072620913
$data[0,3] = 645.8400645940834
$data[1,0] = 45095.99999999999
$data[1,1] = 161
$data[1,2] = -343.9854182427919
$data[1,3] = 680.82845693566
$data[2,0] = 45109.99999999999
$data[2,1] = 170
$data[2,2] = -333.9601283841666
$data[2,3] = 742.8070589525896
$data[3,0] = 45116.99999999999
$data[3,1] = 174
$data[3,2] = -318.6749856658634
$data[3,3] = 687.5062891058373
$data[4,0] = 45123.99999999999
$data[4,1] = 178
$data[4,2] = -352.3714719301074
$data[4,3] = 722.2444814216682
$data[5,0] = 45130.99999999999
$data[5,1] = 183
$data[5,2] = -330.6496146308234
$data[5,3] = 652.3106223543155
$data[6,0] = 45144.99999999999
$data[6,1] = 191
$data[6,2] = -353.7634287715082
$data[6,3] = 730.1491499144286
$data[7,0] = 45151.99999999999
$data[7,1] = 196
$data[7,2] = -321.534885388194
$data[7,3] = 706.0642154769115
$data[8,0] = 45158.99999999999
$data[8,1] = 200
$data[8,2] = -281.0140441386047
$data[8,3] = 712.7719157912479
$data[9,0] = 45165.99999999999
$data[9,1] = 204
$data[9,2] = -294.8176825104684
$data[9,3] = 739.7784785672939
$data[10,0] = 45172.99999999999
$data[10,1] = 209
$data[10,2] = -322.9409993698465
$data[10,3] = 732.6823035914514
$data[11,0] = 45179.99999999999
$data[11,1] = 213
$data[11,2] = -301.4120945782757
$data[11,3] = 728.1928174364883
$data[12,0] = 45186.99999999999
$data[12,1] = 217
$data[12,2] = -263.0442206469735
$data[12,3] = 760.9798970493381
$data[13,0] = 45200.99999999999
$data[13,1] = 226
$data[13,2] = -281.350011021861
$data[13,3] = 719.8597377354893
$data[14,0] = 45207.99999999999
$data[14,1] = 230
$data[14,2] = -268.5273267578253
$data[14,3] = 710.621013820412
$data[15,0] = 45214.99999999999
$data[15,1] = 235
$data[15,2] = -285.9924399756441
$data[15,3] = 788.6777926461093
$data[16,0] = 45221.99999999999
$data[16,1] = 239
$data[16,2] = -290.6420370938
$data[16,3] = 796.8602989310676
$data[17,0] = 45228.99999999999
$data[17,1] = 243
$data[17,2] = -317.6910386993482
$data[17,3] = 780.7931033295355
$data[18,0] = 45235.99999999999
$data[18,1] = 248
$data[18,2] = -250.6482668899592
$data[18,3] = 779.1088612859695
$data[19,0] = 45242.99999999999
$data[19,1] = 252
$data[19,2] = -266.1761866102591
$data[19,3] = 798.927624930774
$data[20,0] = 45249.99999999999
$data[20,1] = 256
$data[20,2] = -232.1623500529986
$data[20,3] = 756.5052600041496
$data[21,0] = 45270.99999999999
$data[21,1] = 269
$data[21,2] = -233.532472083977
$data[21,3] = 782.7690069013039
$data[22,0] = 45298.99999999999
$data[22,1] = 287
$data[22,2] = -251.4564423758537
$data[22,3] = 842.3494371738827
$data[23,0] = 45305.99999999999
$data[23,1] = 291
$data[23,2] = -212.8406281162591
$data[23,3] = 821.6291301310275
$data[24,0] = 45312.99999999999
$data[24,1] = 295
$data[24,2] = -231.8016702062841
$data[24,3] = 843.8777962376334
$data[25,0] = 45319.99999999999
$data[25,1] = 300
$data[25,2] = -220.2981786375785
$data[25,3] = 792.4139293110925
$data[26,0] = 45326.99999999999
$data[26,1] = 304
$data[26,2] = -226.5764976003002
$data[26,3] = 820.7468136687845
$data[27,0] = 45333.99999999999
$data[27,1] = 308
$data[27,2] = -204.2876149493637
$data[27,3] = 833.4721080047974
$data[28,0] = 45375.99999999999
$data[28,1] = 334
$data[28,2] = -153.3741988184265
$data[28,3] = 876.6952235747349
$data[29,0] = 45382.99999999999
$data[29,1] = 339
$data[29,2] = -222.9291419510829
$data[29,3] = 846.3227036657568
$data[30,0] = 45389.99999999999
$data[30,1] = 343
$data[30,2] = -159.697174859825
$data[30,3] = 895.2059385815541
$data[31,0] = 45403.99999999999
$data[31,1] = 352
$data[31,2] = -170.2246999748005
$data[31,3] = 822.5816820080929
$data[32,0] = 45417.99999999999
$data[32,1] = 360
$data[32,2] = -150.5677251278682
$data[32,3] = 898.729242971763
$data[33,0] = 45424.99999999999
$data[33,1] = 365
$data[33,2] = -177.4025284065671
$data[33,3] = 876.4238327926166
$data[34,0] = 45431.99999999999
$data[34,1] = 369
$data[34,2] = -137.8725666660278
$data[34,3] = 905.0119463383094
$data[35,0] = 45438.99999999999
$data[35,1] = 373
$data[35,2] = -148.6853230176696
$data[35,3] = 918.231087655628
$data[36,0] = 45445.99999999999
$data[36,1] = 378
$data[36,2] = -151.0888064892128
$data[36,3] = 913.9654649350545
$data[37,0] = 45452.99999999999
$data[37,1] = 382
$data[37,2] = -154.7618778415726
$data[37,3] = 937.4679860374006
$data[38,0] = 45459.99999999999
$data[38,1] = 386
$data[38,2] = -112.1490089208133
$data[38,3] = 980.5764695270449
$data[39,0] = 45466.99999999999
$data[39,1] = 391
$data[39,2] = -102.2031569977965
$data[39,3] = 888.3715760419113
$data[40,0] = 45473.99999999999
$data[40,1] = 395
$data[40,2] = -135.3332787949707
$data[40,3] = 924.6913492192172
$data[41,0] = 45480.99999999999
$data[41,1] = 399
$data[41,2] = -135.7493431936953
$data[41,3] = 954.6853816803887
$data[42,0] = 45529.99999999999
$data[42,1] = 430
$data[42,2] = -87.76405667224924
$data[42,3] = 937.5789685222209
$data[43,0] = 45536.99999999999
$data[43,1] = 434
$data[43,2] = -80.75062175320909
$data[43,3] = 927.9239801484389
$data[44,0] = 45543.99999999999
$data[44,1] = 438
$data[44,2] = -56.93394194246566
$data[44,3] = 950.1587631493961
$data[45,0] = 45557.99999999999
$data[45,1] = 447
$data[45,2] = -86.43123512672977
$data[45,3] = 1010.082449828148
$data[46,0] = 45564.99999999999
$data[46,1] = 451
$data[46,2] = -75.31943486428588
$data[46,3] = 966.4635371368332
$data[47,0] = 45571.99999999999
$data[47,1] = 456
$data[47,2] = -34.32698890110451
$data[47,3] = 981.4689957172174
$data[48,0] = 45578.99999999999
$data[48,1] = 460
$data[48,2] = -55.5704493853756
$data[48,3] = 951.0253563390744
$data[49,0] = 45585.99999999999
$data[49,1] = 464
$data[49,2] = -68.11992147884322
$data[49,3] = 991.9917427706355
$data[50,0] = 45592.99999999999
$data[50,1] = 469
$data[50,2] = -24.64248984397513
$data[50,3] = 999.7647493371828
$data[51,0] = 45599.99999999999
$data[51,1] = 473
$data[51,2] = -43.34717584680649
$data[51,3] = 949.4572638042836
$data[52,0] = 45606.99999999999
$data[52,1] = 477
$data[52,2] = -43.58668698953687
$data[52,3] = 993.5152087693371
$data[53,0] = 45613.99999999999
$data[53,1] = 482
$data[53,2] = -43.47019230756953
$data[53,3] = 979.4885386447816
$data[54,0] = 45620.99999999999
$data[54,1] = 486
$data[54,2] = -26.14006226312189
$data[54,3] = 997.3961479447844
$data[55,0] = 45627.99999999999
$data[55,1] = 490
$data[55,2] = -23.74471421805152
$data[55,3] = 984.908406748633
$data[56,0] = 45634.99999999999
$data[56,1] = 495
$data[56,2] = -4.169782577778082
$data[56,3] = 1008.623926228689
$data[57,0] = 45641.99999999999
$data[57,1] = 499
$data[57,2] = -42.99950280055904
$data[57,3] = 1004.827648872626
$data[58,0] = 45648.99999999999
$data[58,1] = 503
$data[58,2] = -20.06370552595211
$data[58,3] = 1023.930364516322
$data[59,0] = 45655.99999999999
$data[59,1] = 508
$data[59,2] = -39.30104356882455
$data[59,3] = 1015.72373861081

$wsForecast.Range("A2:D61").Value = $data

# Match the date-column formatting used for column A on the other sheets.
$wsWeekly.Range("A2").Copy()
$wsForecast.Range("A2:A61").PasteSpecial(-4122)

$excel.CutCopyMode = 0
